# Update the "decision reminder" instruction text (F8/F9) so the new-item
# key changes from K to J, and move the selection to F8 — matching the
# commit "instructions and response options changed for recognition".
#
# F10/F11 ("Még 2/1 kérdés van hátra.") keep their existing text; editing
# F8/F9's shared string causes Excel to re-slot the shared-string table so
# those two strings end up at new indices automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$decisionText = @"
A döntésre 4 másodperce lesz.
Minden képet nézzen meg figyelmesen, és minden képre adjon választ, akkor is, ha a döntés nehéz.

A döntését így jelölje:
Régi - F
Új – J
"@

$ws.Cells.Item(8, 6).Value = $decisionText
$ws.Cells.Item(9, 6).Value = $decisionText

$ws.Range("F8").Select()
